$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    8   = -8.981200000000003
    10  = -8.317899999999996
    12  = -7.230199999999996
    18  = -8.850499999999997
    25  = -8.298899999999996
    37  = -7.982000000000003
    55  = -9.030199999999995
    68  = -6.904699999999995
    77  = -6.205299999999995
    78  = -7.799200000000002
    79  = -6.112499999999995
    80  = -7.597899999999997
    81  = -7.471599999999998
    82  = -8.410399999999992
    84  = -8.765800000000004
    101 = -7.958099999999995
    102 = -7.9019
}

foreach ($row in $updates.Keys) {
    $ws.Range("D$row").Value = $updates[$row]
}
